$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44222
$ws.Cells.Item(3, 4).Value = 44222
$ws.Cells.Item(6, 4).Value = 44280
$ws.Cells.Item(7, 4).Value = 44280
$ws.Cells.Item(8, 4).Value = 44272
$ws.Cells.Item(8, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(9, 4).Value = 44272
$ws.Cells.Item(9, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(10, 4).Value = 44350
$ws.Cells.Item(11, 4).Value = 44350
$ws.Cells.Item(12, 4).Value = 44217
$ws.Cells.Item(13, 4).Value = 44217
$ws.Cells.Item(14, 4).Value = 44362
$ws.Cells.Item(15, 4).Value = 44362
$ws.Cells.Item(16, 4).Value = 44405
$ws.Cells.Item(17, 4).Value = 44405
$ws.Cells.Item(18, 4).Value = 44344
$ws.Cells.Item(18, 14).Value = "`$/docena de 1 kilo"
$ws.Cells.Item(19, 4).Value = 44344
$ws.Cells.Item(19, 14).Value = "`$/docena de 1 kilo"
$ws.Cells.Item(20, 4).Value = 44320
$ws.Cells.Item(21, 4).Value = 44320
$ws.Cells.Item(22, 4).Value = 44579
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 4).Value = 44579
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(24, 4).Value = 44285
$ws.Cells.Item(25, 4).Value = 44285
$ws.Cells.Item(26, 4).Value = 44308
$ws.Cells.Item(27, 4).Value = 44308
$ws.Cells.Item(28, 4).Value = 44187
$ws.Cells.Item(29, 4).Value = 44187
$ws.Cells.Item(30, 4).Value = 44160
$ws.Cells.Item(31, 4).Value = 44160
$ws.Cells.Item(32, 4).Value = 44609
$ws.Cells.Item(33, 4).Value = 44609
$ws.Cells.Item(34, 4).Value = 44224
$ws.Cells.Item(35, 4).Value = 44224
$ws.Cells.Item(36, 4).Value = 44400
$ws.Cells.Item(37, 4).Value = 44400
$ws.Cells.Item(38, 4).Value = 44433
$ws.Cells.Item(39, 4).Value = 44433
$ws.Cells.Item(40, 4).Value = 44398
$ws.Cells.Item(40, 15).Value = "Región de Ñuble"
$ws.Cells.Item(41, 4).Value = 44398
$ws.Cells.Item(41, 15).Value = "Región de Ñuble"
$ws.Cells.Item(42, 4).Value = 44194
$ws.Cells.Item(43, 4).Value = 44194
$ws.Cells.Item(44, 4).Value = 44574
$ws.Cells.Item(45, 4).Value = 44574
$ws.Cells.Item(46, 4).Value = 44327
$ws.Cells.Item(47, 4).Value = 44327
$ws.Cells.Item(48, 4).Value = 44306
$ws.Cells.Item(49, 4).Value = 44306
$ws.Cells.Item(50, 4).Value = 44391
$ws.Cells.Item(51, 4).Value = 44391
$ws.Cells.Item(52, 4).Value = 44512
$ws.Cells.Item(52, 10).Value = 200
$ws.Cells.Item(53, 4).Value = 44512
$ws.Cells.Item(53, 10).Value = 100
$ws.Cells.Item(54, 4).Value = 44274
$ws.Cells.Item(55, 4).Value = 44274
$ws.Cells.Item(56, 4).Value = 44442
$ws.Cells.Item(56, 10).Value = 300
$ws.Cells.Item(57, 4).Value = 44442
$ws.Cells.Item(57, 10).Value = 150
$ws.Cells.Item(58, 4).Value = 44292
$ws.Cells.Item(59, 4).Value = 44292
$ws.Cells.Item(60, 4).Value = 44435
$ws.Cells.Item(61, 4).Value = 44435
$ws.Cells.Item(62, 4).Value = 44203
$ws.Cells.Item(62, 10).Value = 200
$ws.Cells.Item(63, 4).Value = 44203
$ws.Cells.Item(63, 10).Value = 100
$ws.Cells.Item(64, 4).Value = 44166
$ws.Cells.Item(65, 4).Value = 44166
$ws.Cells.Item(66, 4).Value = 44355
$ws.Cells.Item(67, 4).Value = 44355
$ws.Cells.Item(68, 4).Value = 44616
$ws.Cells.Item(69, 4).Value = 44616
$ws.Cells.Item(70, 4).Value = 44460
$ws.Cells.Item(71, 4).Value = 44460
$ws.Cells.Item(72, 4).Value = 44341
$ws.Cells.Item(73, 4).Value = 44341
$ws.Cells.Item(74, 4).Value = 44386
$ws.Cells.Item(75, 4).Value = 44386
$ws.Cells.Item(76, 4).Value = 44525
$ws.Cells.Item(77, 4).Value = 44525
$ws.Cells.Item(78, 4).Value = 44316
$ws.Cells.Item(79, 4).Value = 44316
$ws.Cells.Item(80, 4).Value = 44553
$ws.Cells.Item(81, 4).Value = 44553
$ws.Cells.Item(82, 4).Value = 44447
$ws.Cells.Item(83, 4).Value = 44447
$ws.Cells.Item(84, 4).Value = 44237
$ws.Cells.Item(85, 4).Value = 44237
$ws.Cells.Item(86, 4).Value = 44336
$ws.Cells.Item(87, 4).Value = 44336
$ws.Cells.Item(88, 4).Value = 44299
$ws.Cells.Item(89, 4).Value = 44299
$ws.Cells.Item(90, 4).Value = 44239
$ws.Cells.Item(91, 4).Value = 44239
$ws.Cells.Item(92, 4).Value = 44365
$ws.Cells.Item(93, 4).Value = 44365
$ws.Cells.Item(94, 4).Value = 44330
$ws.Cells.Item(94, 15).Value = "Región de Ñuble"
$ws.Cells.Item(95, 4).Value = 44330
$ws.Cells.Item(95, 15).Value = "Región de Ñuble"
$ws.Cells.Item(96, 4).Value = 44490
$ws.Cells.Item(97, 4).Value = 44490
$ws.Cells.Item(98, 4).Value = 44469
$ws.Cells.Item(99, 4).Value = 44469
$ws.Cells.Item(102, 4).Value = 44476
$ws.Cells.Item(103, 4).Value = 44476
$ws.Cells.Item(104, 4).Value = 44425
$ws.Cells.Item(105, 4).Value = 44425
$ws.Cells.Item(106, 4).Value = 44475
$ws.Cells.Item(106, 15).Value = "Región de Ñuble"
$ws.Cells.Item(107, 4).Value = 44475
$ws.Cells.Item(107, 15).Value = "Región de Ñuble"
$ws.Cells.Item(108, 4).Value = 44278
$ws.Cells.Item(108, 10).Value = 300
$ws.Cells.Item(109, 4).Value = 44278
$ws.Cells.Item(109, 10).Value = 150
$ws.Cells.Item(110, 4).Value = 44523
$ws.Cells.Item(111, 4).Value = 44523
$ws.Cells.Item(112, 4).Value = 44231
$ws.Cells.Item(113, 4).Value = 44231
$ws.Cells.Item(114, 4).Value = 44565
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(115, 4).Value = 44565
$ws.Cells.Item(115, 15).Value = "Región Metropolitana"
$ws.Cells.Item(116, 4).Value = 44313
$ws.Cells.Item(117, 4).Value = 44313
$ws.Cells.Item(118, 4).Value = 44453
$ws.Cells.Item(119, 4).Value = 44453
$ws.Cells.Item(120, 4).Value = 44250
$ws.Cells.Item(120, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(121, 4).Value = 44250
$ws.Cells.Item(121, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 4).Value = 44607
$ws.Cells.Item(123, 4).Value = 44607
$ws.Cells.Item(124, 4).Value = 44168
$ws.Cells.Item(125, 4).Value = 44168
$ws.Cells.Item(126, 4).Value = 44582
$ws.Cells.Item(127, 4).Value = 44582
$ws.Cells.Item(128, 4).Value = 44349
$ws.Cells.Item(128, 15).Value = "Región Metropolitana"
$ws.Cells.Item(129, 4).Value = 44349
$ws.Cells.Item(129, 15).Value = "Región Metropolitana"
$ws.Cells.Item(130, 4).Value = 44383
$ws.Cells.Item(130, 15).Value = "Región de Ñuble"
$ws.Cells.Item(131, 4).Value = 44383
$ws.Cells.Item(131, 15).Value = "Región de Ñuble"
$ws.Cells.Item(132, 4).Value = 44334
$ws.Cells.Item(133, 4).Value = 44334
$ws.Cells.Item(134, 4).Value = 44358
$ws.Cells.Item(135, 4).Value = 44358
